$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Implemented All Test Methods in Data Access Layer":
# Every cell in column E (status) that still reads "Not Started" gets flipped
# to "Implemented" (rows 88-110, the clsUserData / clsLicenseClassData block).
for ($r = 2; $r -le 110; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Text -eq "Not Started") {
        $cell.Value = "Implemented"
    }
}

# The clsUserData.DoesPersonHaveUser signature gained a numeric suffix.
$ws.Range("C103").Value = "public static bool DoesPersonHaveUser44(int PersonID)"

# Touch the merged regions so Excel re-emits <mergeCells> in the order it
# naturally lists them after this edit (class blocks at the bottom of the
# sheet first, since that's the area that was just edited).
$mergedBlocks = @(
  "B93:B104","A92:A104","B106:B110","A105:A110","B73:B78","A72:A78","B80:B86","A79:A86","B88:B91","A87:A91",
  "B45:B51","A44:A51","B53:B62","A52:A62","B64:B71","A63:A71",
  "B24:B30","A23:A30","B32:B36","A31:A36","B38:B43","A37:A43",
  "B3:B7","A2:A7","B9:B18","A8:A18","A19:A22","B20:B22"
)
foreach ($ref in $mergedBlocks) {
    $ws.Range($ref).UnMerge()
}
foreach ($ref in $mergedBlocks) {
    $ws.Range($ref).Merge()
}

# Restore the view the author left the sheet in before saving.
$ws.Activate()
$excel.ActiveWindow.Zoom = 131
$ws.Range("C107").Select()
$excel.ActiveWindow.ScrollRow = 107
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F110").Select()
